# Chiffres COVID-19 Valais.xlsx - "Add files via upload" update
# Updates a handful of daily COVID figures (rows 366-373 of the data
# table) and extends the table with two more days of data (rows 372-373
# were previously blank placeholders), then moves the current selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 366: one more "Nb nouveaux décès extra-hospitaliers" (M) ---
# Column M uses a Text ("@") number format, so a plain .Value assignment
# would store the number as text. Flip to General while writing the
# value, then restore the Text format so the style index is unchanged.
$ws.Range("M366").NumberFormat = "General"
$ws.Range("M366").Value = 1
$ws.Range("M366").NumberFormat = "@"

# --- Row 370: corrected "Cumul cas positifs" delta (C) ---
$ws.Range("C370").Value = 88

# --- Row 371: corrected "Cumul cas positifs" delta (C) ---
$ws.Range("C371").Value = 63

# --- Row 372: new day of data (2020-03-12 row in the sheet -> A372=44258) ---
$ws.Range("C372").Value = 43
$ws.Range("E372").Value = 8
$ws.Range("F372").Value = 5
$ws.Range("G372").Value = 26
$ws.Range("L372").Value = "1"
$ws.Range("M372").NumberFormat = "General"
$ws.Range("M372").Value = 0
$ws.Range("M372").NumberFormat = "@"

# --- Row 373: new day of data (A373=44259) ---
$ws.Range("C373").Value = 5
$ws.Range("E373").Value = 9
$ws.Range("F373").Value = 7
$ws.Range("G373").Value = 29
$ws.Range("L373").Value = "0"
$ws.Range("M373").Value = "0"

# --- Move the active selection in the frozen (bottom-right) pane ---
$ws.Range("S377").Select()
